$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: fully replace a paragraph's text (excluding the paragraph
# mark) using an explicit numeric Range - Paragraph.Range.Text = ...
# only replaces the first run, so we avoid it.
# ------------------------------------------------------------------
function Set-ParaText($paraIndex, $text) {
    $p = $d.Paragraphs($paraIndex)
    $start = $p.Range.Start
    $end = $p.Range.End
    $r = $d.Range($start, $end - 1)
    $r.Text = $text
}

# ====================================================================
# 1) Paragraph 3 ("Stopped: ...") : replace text, then append 5 new
#    sub-bullets (ilvl=1, numId=3) describing the issues found.
# ====================================================================
Set-ParaText 3 "Stopped: script_data, the call to get_create_table_from_sys_tables:"

$p3 = $d.Paragraphs(3)
$anchor = $p3.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)
$p4.Range.ListFormat.ListLevelNumber = 2
$s = $p4.Range.Start
$e = $p4.Range.End
$d.Range($s, $e - 1).Text = "Doesn’t put comma after each column"

$p4b = $d.Paragraphs(4)
$p4b.Range.Collapse(0)
$p4b.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(5)
$p5.Range.ListFormat.ListLevelNumber = 2
$s = $p5.Range.Start
$e = $p5.Range.End
$d.Range($s, $e - 1).Text = "It’s the classic table name, not the override temp table name we gave it"

$p5b = $d.Paragraphs(5)
$p5b.Range.Collapse(0)
$p5b.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs(6)
$p6.Range.ListFormat.ListLevelNumber = 2
$s = $p6.Range.Start
$e = $p6.Range.End
$d.Range($s, $e - 1).Text = "It gives  afull create table, including create index… not needed here. See old code"

$p6b = $d.Paragraphs(6)
$p6b.Range.Collapse(0)
$p6b.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs(7)
$p7.Range.ListFormat.ListLevelNumber = 2
$s = $p7.Range.Start
$e = $p7.Range.End
$d.Range($s, $e - 1).Text = "The ALTER TABLE after it: again, use new, data table name"

$p7b = $d.Paragraphs(7)
$p7b.Range.Collapse(0)
$p7b.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)
$p8.Range.ListFormat.ListLevelNumber = 2
$s = $p8.Range.Start
$e = $p8.Range.End
$d.Range($s, $e - 1).Text = "SET _CmprState_=RowState.EXTRA1  - should set the value of the enum"

# ====================================================================
# 2) Paragraph 9 ("Then: ...") : "inserts and updates" -> "updates,
#    deletes"; then append a new bullet (ilvl=0, numId=3) right after.
# ====================================================================
Set-ParaText 9 "Then: when its working, we got the code for insertions. Is it running? Huge. Then give it to claude as guideline for next round: updates, deletes"

$p9 = $d.Paragraphs(9)
$p9.Range.Collapse(0)
$p9.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs(10)
$p10.Range.ListFormat.ListLevelNumber = 1
$s = $p10.Range.Start
$e = $p10.Range.End
$d.Range($s, $e - 1).Text = "Test the whole ‘data table that was empty’ , bulk-insert into it"

# ====================================================================
# 3) After "Smallie: title of server.dbname ..." add a new bullet
#    (ilvl=0, numId=10) about transactional execution.
# ====================================================================
$p20 = $d.Paragraphs(20)
$p20.Range.Collapse(0)
$p20.Range.InsertParagraphAfter()
$p21 = $d.Paragraphs(21)
$p21.Range.ListFormat.ListLevelNumber = 1
$s = $p21.Range.Start
$e = $p21.Range.End
$d.Range($s, $e - 1).Text = "Execution as transaction, so if one fails, it all fails. (by default no)"
